# website & logboek update
# Add a new logboek entry (row 19) for 01-05-2017: "13.00 / School / werken aan de website"
# and update the active sheet view (scroll position + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 19 -----------------------------------------------------------
# A19: date 01-05-2017, formatted the same way as the other date cells (copy
# the format from A17 so the existing date style is reused instead of a new
# one being minted), then assign just the date (no time-of-day component).
$ws.Range("A17").Copy()
$ws.Range("A19").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A19").Value = (Get-Date -Year 2017 -Month 5 -Day 1).Date

# C19: "13.00" — copy the existing "13.00" text cell (C11) so it reuses the
# same shared-string entry and keeps the default (unstyled) cell format.
$ws.Range("C11").Copy()
$ws.Range("C19").PasteSpecial(-4104)  # xlPasteAll

# E19: "School" — copy from E17 (same shared string as every other row).
$ws.Range("E17").Copy()
$ws.Range("E19").PasteSpecial(-4104)  # xlPasteAll

# G19: new remark text.
$ws.Range("G19").Value = "werken aan de website"

# --- Sheet view -------------------------------------------------------------
try { $excel.ActiveWindow.TopLeftCell = $ws.Range("A10") } catch {}
try { $excel.ActiveWindow.ScrollRow = 10 } catch {}
$ws.Range("L19").Select()
